$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- Append the new log row (row 45) ---
$reply = @"
Beste klant,
Bedankt voor uw e-mail. Mijn excuses dat uw bestelling incompleet is geleverd. Om dit voor u op te lossen, heb ik wat meer informatie nodig. Kunt u mij alstublieft het volgende verstrekken:
- Uw bestelnummer
- De ontbrekende item(s)
Met deze gegevens kan ik direct voor u nakijken wat er mis is gegaan en een passende oplossing bieden.
Ik hoor graag van u.
Met vriendelijke groet,
[Naam]
E-mailassistent
"@

$ws.Range("A45").Value = "Klacht over levering"
$ws.Range("B45").Value = "mailmind.test@zohomail.eu"
$ws.Range("C45").Value = "Mijn bestelling is incompleet geleverd. Graag hoor ik hoe dit wordt opgelost."
$ws.Range("D45").Value = "Bestelling / Levering"
$ws.Range("E45").Value = $reply
$ws.Range("F45").Value = "2025-06-22 21:43:42"
$ws.Range("G45").Value = "Ja"

# --- Extend conditional formatting ranges to include the new row ---
$fcsD = $ws.Range("D2:D44").FormatConditions
for ($i = 1; $i -le $fcsD.Count; $i++) {
    $fcsD.Item($i).ModifyAppliesToRange($ws.Range("D2:D45"))
}

$fcsG = $ws.Range("G2:G44").FormatConditions
for ($i = 1; $i -le $fcsG.Count; $i++) {
    $fcsG.Item($i).ModifyAppliesToRange($ws.Range("G2:G45"))
}

# --- Update the Dashboard category counts to reflect the new row ---
# Before: row12=Juridisch/Contract(2), row13=Klacht/Probleem(1),
#         row14=Uitnodiging/Evenement(1), row15=Bestelling/Levering(1)
# After:  row12=Bestelling/Levering(2), row13=Juridisch/Contract(2),
#         row14=Uitnodiging/Evenement(1), row15=Klacht/Probleem(1)
$dash.Range("A12").Value = "Bestelling / Levering"
$dash.Range("B12").Value = 2

$dash.Range("A13").Value = "Juridisch / Contract"
$dash.Range("B13").Value = 2

$dash.Range("A15").Value = "Klacht / Probleem"
$dash.Range("B15").Value = 1
